# Reproduce the "add common excel diff" edit: populate Sheet3 (the only
# sheet in the workbook) with a handful of sample values/strings and
# autofit the columns that hold them, then leave the selection on H16
# with the view scrolled so row 7 is at the top (as in the target file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new shared strings in the same order they appear in the
# target xl/sharedStrings.xml (0: 阿斯顿撒打算, 1: 嘻嘻嘻, 2: xx) by
# setting the cells that reference them in that order first.
$ws.Range("J11").Value = "阿斯顿撒打算"
$ws.Range("F16").Value = "嘻嘻嘻"
$ws.Range("B10").Value = "xx"

# Numeric cells.
$ws.Range("B4").Value = 1111111111
$ws.Range("C10").Value = 222222222
$ws.Range("F11").Value = 333333333333

# Columns B, C and F need to grow to fit their new contents.
$ws.Columns("B").AutoFit()
$ws.Columns("C").AutoFit()
$ws.Columns("F").AutoFit()

# Scroll so row 7 is the first visible row, and finish with H16 selected
# (matches topLeftCell="A7" / selection activeCell="H16" sqref="H16").
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("H16").Select() | Out-Null
